$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, pushing existing rows 66-76 down to 67-77
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new data record
$ws.Cells.Item(66, 1).Value = 4
$ws.Cells.Item(66, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(66, 3).Value = "Los Lagos"
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 4).Value = 44491
$ws.Cells.Item(66, 5).Value = 10
$ws.Cells.Item(66, 6).Value = 100112022
$ws.Cells.Item(66, 7).Value = "Arveja Verde"
$ws.Cells.Item(66, 8).Value = "Perfection"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 80
$ws.Cells.Item(66, 11).Value = 20000
$ws.Cells.Item(66, 12).Value = 20000
$ws.Cells.Item(66, 13).Value = 20000
$ws.Cells.Item(66, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(66, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(66, 16).Value = 800
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"
